# Auto-generated script applying updated market-board price/profit data
# to the Leve profit tracker sheets (columns H-N), per scheduled runner diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 92.77778000000001
$ws.Range("I9").Value = 88.75
$ws.Range("K9").Value = 88.75
$ws.Range("M9").Value = 80.25
$ws.Range("H15").Value = 11401.438
$ws.Range("I15").Value = 11401.438
$ws.Range("K15").Value = 34204.314
$ws.Range("M15").Value = -34035.314
$ws.Range("H28").Value = 1152.0834
$ws.Range("I28").Value = 238.47058
$ws.Range("J28").Value = 3370.8572
$ws.Range("K28").Value = 238.47058
$ws.Range("L28").Value = 3370.8572
$ws.Range("M28").Value = 246.52942
$ws.Range("N28").Value = -4340.8572
$ws.Range("H33").Value = 9524066
$ws.Range("I33").Value = 12345951
$ws.Range("J33").Value = 204.875
$ws.Range("K33").Value = 12345951
$ws.Range("L33").Value = 204.875
$ws.Range("M33").Value = -12345722
$ws.Range("N33").Value = -662.875
$ws.Range("H62").Value = 1084.8
$ws.Range("I62").Value = 1750
$ws.Range("J62").Value = 641.3333
$ws.Range("K62").Value = 1750
$ws.Range("L62").Value = 641.3333
$ws.Range("M62").Value = -1126
$ws.Range("N62").Value = -1889.3333
$ws.Range("H65").Value = 1084.8
$ws.Range("I65").Value = 1750
$ws.Range("J65").Value = 641.3333
$ws.Range("K65").Value = 8750
$ws.Range("L65").Value = 3206.6665
$ws.Range("M65").Value = -5630
$ws.Range("N65").Value = -9446.666499999999
$ws.Range("H111").Value = 4788.1665
$ws.Range("I111").Value = 5182.25
$ws.Range("J111").Value = 4000
$ws.Range("K111").Value = 15546.75
$ws.Range("L111").Value = 12000
$ws.Range("M111").Value = -12479.75
$ws.Range("N111").Value = -18134
$ws.Range("H116").Value = 3571
$ws.Range("I116").Value = 799.4
$ws.Range("K116").Value = 799.4
$ws.Range("M116").Value = 2642.6
$ws.Range("H132").Value = 2917471.2
$ws.Range("I132").Value = 3761561.8
$ws.Range("J132").Value = 1521.8182
$ws.Range("K132").Value = 11284685.4
$ws.Range("L132").Value = 4565.4546
$ws.Range("M132").Value = -11282155.4
$ws.Range("N132").Value = -9625.454600000001
$ws.Range("H135").Value = 3604.3333
$ws.Range("I135").Value = 644.7742
$ws.Range("K135").Value = 5802.967799999999
$ws.Range("M135").Value = -3267.967799999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17396.965
$ws.Range("I32").Value = 17157.688
$ws.Range("J32").Value = 19700
$ws.Range("K32").Value = 17157.688
$ws.Range("L32").Value = 19700
$ws.Range("M32").Value = -16870.688
$ws.Range("N32").Value = -20274
$ws.Range("H45").Value = 83334590
$ws.Range("I45").Value = 83334590
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 83334590
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -83334213
$ws.Range("H61").Value = 1588.2122
$ws.Range("I61").Value = 1238.8846
$ws.Range("K61").Value = 1238.8846
$ws.Range("M61").Value = -1026.8846
$ws.Range("H63").Value = 2500651.2
$ws.Range("I63").Value = 2500651.2
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2500651.2
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -2499965.2
$ws.Range("H66").Value = 2500651.2
$ws.Range("I66").Value = 2500651.2
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 12503256
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -12499824
$ws.Range("H74").Value = 787.4545000000001
$ws.Range("I74").Value = 526.9091
$ws.Range("J74").Value = 1048
$ws.Range("K74").Value = 526.9091
$ws.Range("L74").Value = 1048
$ws.Range("M74").Value = 347.0909
$ws.Range("N74").Value = -2796
$ws.Range("H77").Value = 787.4545000000001
$ws.Range("I77").Value = 526.9091
$ws.Range("J77").Value = 1048
$ws.Range("K77").Value = 2634.5455
$ws.Range("L77").Value = 5240
$ws.Range("M77").Value = 1733.4545
$ws.Range("N77").Value = -13976
$ws.Range("H132").Value = 3470.238
$ws.Range("I132").Value = 3641.102
$ws.Range("J132").Value = 2872.2144
$ws.Range("K132").Value = 10923.306
$ws.Range("L132").Value = 8616.643199999999
$ws.Range("M132").Value = -8393.306
$ws.Range("N132").Value = -13676.6432
$ws.Range("H136").Value = 1588.2122
$ws.Range("I136").Value = 1238.8846
$ws.Range("K136").Value = 3716.6538
$ws.Range("M136").Value = -1166.6538

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2316.3333
$ws.Range("I99").Value = 2379.6
$ws.Range("K99").Value = 2379.6
$ws.Range("M99").Value = -881.5999999999999
$ws.Range("H134").Value = 2013.1666
$ws.Range("I134").Value = 1805.6538
$ws.Range("J134").Value = 2258.4092
$ws.Range("K134").Value = 5416.9614
$ws.Range("L134").Value = 6775.2276
$ws.Range("M134").Value = -2881.9614
$ws.Range("N134").Value = -11845.2276

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2453485.8
$ws.Range("I31").Value = 1799.6274
$ws.Range("K31").Value = 1799.6274
$ws.Range("M31").Value = -1504.6274
$ws.Range("H34").Value = 2453485.8
$ws.Range("I34").Value = 1799.6274
$ws.Range("K34").Value = 1799.6274
$ws.Range("M34").Value = -1597.6274
$ws.Range("H58").Value = 1062.1
$ws.Range("I58").Value = 1118.9286
$ws.Range("J58").Value = 929.5
$ws.Range("K58").Value = 1118.9286
$ws.Range("L58").Value = 929.5
$ws.Range("M58").Value = -915.9286
$ws.Range("N58").Value = -1335.5
$ws.Range("H99").Value = 2762.4285
$ws.Range("I99").Value = 2527.7144
$ws.Range("J99").Value = 2997.1428
$ws.Range("K99").Value = 2527.7144
$ws.Range("L99").Value = 2997.1428
$ws.Range("M99").Value = -1029.7144
$ws.Range("N99").Value = -5993.1428
$ws.Range("H107").Value = 591.3509
$ws.Range("I107").Value = 553.5135
$ws.Range("J107").Value = 661.35
$ws.Range("K107").Value = 553.5135
$ws.Range("L107").Value = 661.35
$ws.Range("M107").Value = 1366.4865
$ws.Range("N107").Value = -4501.35
$ws.Range("H126").Value = 2762.4285
$ws.Range("I126").Value = 2527.7144
$ws.Range("J126").Value = 2997.1428
$ws.Range("K126").Value = 7583.1432
$ws.Range("L126").Value = 8991.428400000001
$ws.Range("M126").Value = -5113.1432
$ws.Range("N126").Value = -13931.4284
$ws.Range("H132").Value = 2405344.5
$ws.Range("I132").Value = 740.4
$ws.Range("J132").Value = 7355999.5
$ws.Range("K132").Value = 2221.2
$ws.Range("L132").Value = 22067998.5
$ws.Range("M132").Value = 308.8000000000002
$ws.Range("N132").Value = -22073058.5
$ws.Range("H134").Value = 630.4
$ws.Range("I134").Value = 635.6842
$ws.Range("J134").Value = 601.7143
$ws.Range("K134").Value = 1907.0526
$ws.Range("L134").Value = 1805.1429
$ws.Range("M134").Value = 627.9474
$ws.Range("N134").Value = -6875.1429
$ws.Range("H136").Value = 1062.1
$ws.Range("I136").Value = 1118.9286
$ws.Range("J136").Value = 929.5
$ws.Range("K136").Value = 3356.7858
$ws.Range("L136").Value = 2788.5
$ws.Range("M136").Value = -806.7857999999997
$ws.Range("N136").Value = -7888.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1870984.2
$ws.Range("I131").Value = 14654.143
$ws.Range("J131").Value = 2116160
$ws.Range("K131").Value = 43962.429
$ws.Range("L131").Value = 6348480
$ws.Range("M131").Value = -38922.429
$ws.Range("N131").Value = -6358560

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 31252460
$ws.Range("I122").Value = 66669530
$ws.Range("J122").Value = 2105.7646
$ws.Range("K122").Value = 200008590
$ws.Range("L122").Value = 6317.293799999999
$ws.Range("M122").Value = -200006140
$ws.Range("N122").Value = -11217.2938
$ws.Range("H126").Value = 953894.7
$ws.Range("I126").Value = 1719.3684
$ws.Range("J126").Value = 2084602.9
$ws.Range("K126").Value = 5158.1052
$ws.Range("L126").Value = 6253808.699999999
$ws.Range("M126").Value = -2688.1052
$ws.Range("N126").Value = -6258748.699999999
$ws.Range("H131").Value = 25500.5
$ws.Range("J131").Value = 25500.5
$ws.Range("L131").Value = 25500.5
$ws.Range("N131").Value = -35580.5
$ws.Range("H132").Value = 23384.979
$ws.Range("I132").Value = 31119.941
$ws.Range("J132").Value = 3155.077
$ws.Range("K132").Value = 93359.823
$ws.Range("L132").Value = 9465.231
$ws.Range("M132").Value = -90829.823
$ws.Range("N132").Value = -14525.231

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 35716610
$ws.Range("I7").Value = 2375
$ws.Range("K7").Value = 2375
$ws.Range("M7").Value = -2263
$ws.Range("H126").Value = 35716610
$ws.Range("I126").Value = 2375
$ws.Range("K126").Value = 7125
$ws.Range("M126").Value = -4655
$ws.Range("H132").Value = 7905.0884
$ws.Range("I132").Value = 18217.5
$ws.Range("K132").Value = 54652.5
$ws.Range("M132").Value = -52122.5
$ws.Range("H133").Value = 25940.46
$ws.Range("J133").Value = 25940.46
$ws.Range("L133").Value = 25940.46
$ws.Range("N133").Value = -31000.46
$ws.Range("H136").Value = 5834.483
$ws.Range("I136").Value = 7894.4116
$ws.Range("J136").Value = 2916.25
$ws.Range("K136").Value = 23683.2348
$ws.Range("L136").Value = 8748.75
$ws.Range("M136").Value = -21133.2348
$ws.Range("N136").Value = -13848.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1345.902
$ws.Range("I132").Value = 1187.0256
$ws.Range("K132").Value = 3561.0768
$ws.Range("M132").Value = -1031.0768
$ws.Range("H136").Value = 1278.4807
$ws.Range("I136").Value = 1452.6666
$ws.Range("J136").Value = 1090.36
$ws.Range("K136").Value = 4357.9998
$ws.Range("L136").Value = 3271.08
$ws.Range("M136").Value = -1807.9998
$ws.Range("N136").Value = -8371.08
